$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 10: ORG_CODE (C10) and TEAM_Code (O10) change from "LT52" to "LT51";
# TEAM_NAME (P10) changes to match the new team code's name.
$ws.Range("C10").Value = "LT51"
$ws.Range("O10").Value = "LT51"
$ws.Range("P10").Value = "LT51 - Services & Call Center Team"

# Update the view state (scroll position / selection) to match the edited file.
$ws.Range("T11").Select()
$ws.Application.ActiveWindow.ScrollColumn = 13
